$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("661-17165",    "AIRPODS PRO 1ST GEN, RIGHT, CI/AR",                "349.411,50"),
    @("Z661-13579",   "SPACE GRAY, 128GB, IPAD 7TH GEN, WI-FI, INT",      "1.009.873,75"),
    @("LL661-42726",  "Display, iPhone 16 Pro, USA",                      "1.380.293,95"),
    @("923-03869",    "AIRPODS PRO 1ST GEN, EAR TIPS, MEDIUM",            "15.462,00"),
    @("923-04957",    "BLACK, AIRPODS MAX, LEFT, EAR CUSHION",            "142.941,60"),
    @("923-03869",    "AIRPODS PRO 1ST GEN, EAR TIPS, MEDIUM",            "15.462,00"),
    @("923-03869",    "AIRPODS PRO 1ST GEN, EAR TIPS, MEDIUM",            "15.462,00"),
    @("923-03869",    "AIRPODS PRO 1ST GEN, EAR TIPS, MEDIUM",            "15.462,00")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 6).Value = $data[$i][0]
    $ws.Cells.Item($row, 7).Value = $data[$i][1]
    $ws.Cells.Item($row, 8).Value = $data[$i][2]
}
